# Split the three long "Programa"/"Bibliografia" paragraphs in LOB1021.docx
# into one run per item, separated by manual line breaks (<w:br/>), matching
# the upstream site-build formatting change.
$d = $word.ActiveDocument

# --- Paragraph 1: Portuguese "Programa" course outline ---
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("Óptica Geométrica: conceitos básicos.", $false, $false, $false, $false, $false, $true, 1, $false)
if (-not $found1) { throw "Portuguese programa paragraph not found" }
$para1 = $rng1.Paragraphs(1).Range
$para1.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">1) Óptica Geométrica: conceitos básicos. </w:t><w:br/><w:t>2) Interferência: a experiência de Young; coerência; figuras de interferência; o interferômetro de Michelson.</w:t><w:br/><w:t>3) Difração.</w:t><w:br/><w:t>4) Polarização.</w:t><w:br/><w:t>5) Relatividade: os postulados da relatividade, as transformações de Lorentz, simultaneidade, tempo e comprimento; momento linear, trabalho e energia;</w:t><w:br/><w:t xml:space="preserve">6) Primórdios da teoria quântica: a hipótese de Plank; o efeito fotoelétrico, quantização do fóton; ondas de De Broglie, o efeito Compton, a difração de elétrons, interferência; </w:t><w:br/><w:t>7) Princípios básicos da mecânica quântica: o princípio de incerteza; a equação de Schrödinger.</w:t></w:r></w:p>')

# --- Paragraph 2: English "Programa" course outline (italic) ---
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Geometrical Optics: basic concepts.", $false, $false, $false, $false, $false, $true, 1, $false)
if (-not $found2) { throw "English programa paragraph not found" }
$para2 = $rng2.Paragraphs(1).Range
$para2.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:i/></w:rPr><w:t>1) Geometrical Optics: basic concepts.</w:t><w:br/><w:t>2) Interference: Young''s experience; coherence; interference figures; the Michelson interferometer.</w:t><w:br/><w:t>3) Diffraction.</w:t><w:br/><w:t>4) Polarization.</w:t><w:br/><w:t>5) Relativity: the postulates of relativity, Lorentz transformations, simultaneity, time and length; linear momentum, work and energy;</w:t><w:br/><w:t>6) Early days of quantum theory: the hypothesis of Planck; the photoelectric effect, quantization of the photon; De Broglie waves, the Compton effect, the electron diffraction, interference;</w:t><w:br/><w:t>7) Basic principles of quantum mechanics: the uncertainty principle; the Schrödinger equation.</w:t></w:r></w:p>')

# --- Paragraph 3: Bibliography ---
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("NUSSENZVEIG, H.M. Curso de Física Básica.", $false, $false, $false, $false, $false, $true, 1, $false)
if (-not $found3) { throw "Bibliography paragraph not found" }
$para3 = $rng3.Paragraphs(1).Range
$para3.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>NUSSENZVEIG, H.M. Curso de Física Básica. Vol. 4, Edgard Blucher (2008).</w:t><w:br/><w:t>RESNICK, R.; HALLIDAY, D. Fundamentos de Física. Vol. 4, LTC (2008).</w:t><w:br/><w:t>TIPLER, P.; MOSCA, G. Física para Cientistas e Engenheiros. Vol. 4, LTC (2008).</w:t><w:br/><w:t>SEARS, F. W.; ZEMANSKY, M. W.; YOUNG, H. D.; FREEDMAN, R. A. Física IV, Vol. 4, Pearson Addison Wesley (2009).</w:t><w:br/><w:t>JEWETT Jr, John W.; SERWAY, Raymond A. Princípios de Física. Vol. 4, Thomson Pioneira (2008).</w:t></w:r></w:p>')

Write-Host "Done"
